$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell that carries the default (unstyled) cell format,
# used to keep forced-text cells on the original style index.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "62.124.02"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "2.417.23"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.85"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.96"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +3.11%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "2.413.29"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.18"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  +4.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000176"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +5.32%  "
$ws.Range("D16").Value = "2.861.37"
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("D17").Value = "61.930.84"
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "2.415.72"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.15"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = "  +2.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.20"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.30"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.43"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.05"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +8.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "590.43"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +16.14%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "2.520.95"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("E30").Value = "  +5.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.29"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("E32").Value = "  +5.29%  "
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("E35").Value = "  +2.79%  "
$ws.Range("E36").Value = "  +5.61%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.78"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +2.64%  "
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.69"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.46"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("E44").Value = "  +12.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "151.21"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("E47").Value = "  +3.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.23"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  +4.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.591"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +2.69%  "
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("E51").Value = "  +2.15%  "
